# Update TPM-derived values in the LR-pairs sheet for Lta-Tnfrsf1b
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.014142
$ws.Range("H2").Value = 0.042426
$ws.Range("M2").Value = 5.916202333333334
$ws.Range("N2").Value = 17.748607
$ws.Range("O2").Value = 0.3515586392055965
$ws.Range("P2").Value = 0.3515586392055965
$ws.Range("Q2").Value = 0.08366693339800001
$ws.Range("R2").Value = 0.753002400582
$ws.Range("S2").Value = 0.3515586392055965
$ws.Range("T2").Value = 0.3515586392055965

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.014142
$ws.Range("H3").Value = 0.042426
$ws.Range("O3").Value = 0.6159539016771971
$ws.Range("P3").Value = 0.6159539016771971
$ws.Range("Q3").Value = 0.146589980506
$ws.Range("R3").Value = 1.319309824554
$ws.Range("S3").Value = 0.6159539016771971
$ws.Range("T3").Value = 0.6159539016771971

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.014142
$ws.Range("H4").Value = 0.042426
$ws.Range("M4").Value = 0.5467150000000001
$ws.Range("O4").Value = 0.03248745911720639
$ws.Range("P4").Value = 0.03248745911720639
$ws.Range("Q4").Value = 0.007731643530000001
$ws.Range("R4").Value = 0.06958479177
$ws.Range("S4").Value = 0.03248745911720639
$ws.Range("T4").Value = 0.03248745911720639
